# Scheduled-runner market-price refresh for the Atomos_Profits Leve sheets.
# Re-pulls currentAveragePrice*/LevePrice*/LeveProfit* (cols H-N) for a handful
# of rows across the eight job sheets (one per DoH crafting job).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1467.1052
$ws.Range("I6").Value = 101.5
$ws.Range("K6").Value = 304.5
$ws.Range("M6").Value = -192.5

# Row 32
$ws.Range("H32").Value = 441.7143
$ws.Range("J32").Value = 498.4
$ws.Range("L32").Value = 498.4
$ws.Range("N32").Value = -1150.4

# Row 87
$ws.Range("H87").Value = 26471.6
$ws.Range("J87").Value = 26471.6
$ws.Range("L87").Value = 26471.6
$ws.Range("N87").Value = -28967.6

# Row 90
$ws.Range("H90").Value = 26471.6
$ws.Range("J90").Value = 26471.6
$ws.Range("L90").Value = 79414.79999999999
$ws.Range("N90").Value = -91894.79999999999

# Row 95
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492

# Row 98
$ws.Range("H98").Value = 2744.7856
$ws.Range("I98").Value = 2821.8096
$ws.Range("J98").Value = 2513.7144
$ws.Range("K98").Value = 2821.8096
$ws.Range("L98").Value = 2513.7144
$ws.Range("M98").Value = -1323.8096
$ws.Range("N98").Value = -5509.7144

# Row 105
$ws.Range("H105").Value = 29900
$ws.Range("J105").Value = 29900
$ws.Range("L105").Value = 29900
$ws.Range("N105").Value = -36888

# Row 112
$ws.Range("H112").Value = 8929936
$ws.Range("J112").Value = 10417713
$ws.Range("L112").Value = 31253139
$ws.Range("N112").Value = -31255355

# Row 122
$ws.Range("H122").Value = 2744.7856
$ws.Range("I122").Value = 2821.8096
$ws.Range("J122").Value = 2513.7144
$ws.Range("K122").Value = 8465.4288
$ws.Range("L122").Value = 7541.1432
$ws.Range("M122").Value = -6015.4288
$ws.Range("N122").Value = -12441.1432

# Row 125
$ws.Range("H125").Value = 1957.4546
$ws.Range("I125").Value = 1557.5
$ws.Range("K125").Value = 14017.5
$ws.Range("M125").Value = -11557.5

# Row 133
$ws.Range("H133").Value = 22401.117
$ws.Range("J133").Value = 22401.117
$ws.Range("L133").Value = 22401.117
$ws.Range("N133").Value = -32521.117

# Row 137
$ws.Range("H137").Value = 2894.5398
$ws.Range("I137").Value = 2940.5227
$ws.Range("J137").Value = 2788.0527
$ws.Range("K137").Value = 8821.5681
$ws.Range("L137").Value = 8364.158100000001
$ws.Range("M137").Value = -6271.5681
$ws.Range("N137").Value = -13464.1581

# Row 138
$ws.Range("H138").Value = 2275.776
$ws.Range("I138").Value = 1509.4375
$ws.Range("J138").Value = 3218.9614
$ws.Range("K138").Value = 4528.3125
$ws.Range("L138").Value = 9656.8842
$ws.Range("M138").Value = 611.6875
$ws.Range("N138").Value = -19936.8842

# Row 141
$ws.Range("H141").Value = 714725.2
$ws.Range("J141").Value = 1426919.1
$ws.Range("L141").Value = 4280757.300000001
$ws.Range("N141").Value = -4291117.300000001

$ws = $wb.Worksheets.Item("ARM")
# Row 29
$ws.Range("H29").Value = 27003.334
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 80010
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 80010
$ws.Range("M29").Value = -192
$ws.Range("N29").Value = -80626

# Row 32
$ws.Range("H32").Value = 4672.2373
$ws.Range("I32").Value = 4555.3936
$ws.Range("J32").Value = 8333.333000000001
$ws.Range("K32").Value = 4555.3936
$ws.Range("L32").Value = 8333.333000000001
$ws.Range("M32").Value = -4268.3936
$ws.Range("N32").Value = -8907.333000000001

# Row 110
$ws.Range("H110").Value = 2432.3333
$ws.Range("I110").Value = 712.25
$ws.Range("K110").Value = 712.25
$ws.Range("M110").Value = 1332.75

# Row 132
$ws.Range("H132").Value = 3662.6765
$ws.Range("I132").Value = 3773.5908
$ws.Range("J132").Value = 3459.3333
$ws.Range("K132").Value = 11320.7724
$ws.Range("L132").Value = 10377.9999
$ws.Range("M132").Value = -8790.7724
$ws.Range("N132").Value = -15437.9999

# Row 134
$ws.Range("H134").Value = 29562.5
$ws.Range("J134").Value = 29562.5
$ws.Range("L134").Value = 29562.5
$ws.Range("N134").Value = -39702.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2362.8462
$ws.Range("I20").Value = 2421.2
$ws.Range("K20").Value = 2421.2
$ws.Range("M20").Value = -2174.2

# Row 64
$ws.Range("H64").Value = 683
$ws.Range("I64").Value = 699.75
$ws.Range("J64").Value = 673.4286
$ws.Range("K64").Value = 699.75
$ws.Range("L64").Value = 673.4286
$ws.Range("M64").Value = -474.75
$ws.Range("N64").Value = -1123.4286

# Row 67
$ws.Range("H67").Value = 683
$ws.Range("I67").Value = 699.75
$ws.Range("J67").Value = 673.4286
$ws.Range("K67").Value = 699.75
$ws.Range("L67").Value = 673.4286
$ws.Range("M67").Value = 80.25
$ws.Range("N67").Value = -2233.4286

# Row 74
$ws.Range("H74").Value = 50973.6
$ws.Range("J74").Value = 50973.6
$ws.Range("L74").Value = 50973.6
$ws.Range("N74").Value = -52845.6

# Row 77
$ws.Range("H77").Value = 50973.6
$ws.Range("J77").Value = 50973.6
$ws.Range("L77").Value = 152920.8
$ws.Range("N77").Value = -162280.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2588.6667
$ws.Range("I16").Value = 1714.1666
$ws.Range("K16").Value = 1714.1666
$ws.Range("M16").Value = -1427.1666

# Row 31
$ws.Range("H31").Value = 6009.6665
$ws.Range("I31").Value = 5267.778
$ws.Range("J31").Value = 6566.0835
$ws.Range("K31").Value = 5267.778
$ws.Range("L31").Value = 6566.0835
$ws.Range("M31").Value = -4972.778
$ws.Range("N31").Value = -7156.0835

# Row 34
$ws.Range("H34").Value = 6009.6665
$ws.Range("I34").Value = 5267.778
$ws.Range("J34").Value = 6566.0835
$ws.Range("K34").Value = 5267.778
$ws.Range("L34").Value = 6566.0835
$ws.Range("M34").Value = -5065.778
$ws.Range("N34").Value = -6970.0835

# Row 113
$ws.Range("H113").Value = 2588.6667
$ws.Range("I113").Value = 1714.1666
$ws.Range("K113").Value = 1714.1666
$ws.Range("M113").Value = 455.8334

# Row 134
$ws.Range("H134").Value = 22045.4
$ws.Range("I134").Value = 29252.285
$ws.Range("J134").Value = 5229.3335
$ws.Range("K134").Value = 87756.855
$ws.Range("L134").Value = 15688.0005
$ws.Range("M134").Value = -85221.855
$ws.Range("N134").Value = -20758.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1249.5
$ws.Range("I5").Value = 479.4
$ws.Range("J5").Value = 5100
$ws.Range("K5").Value = 1438.2
$ws.Range("L5").Value = 15300
$ws.Range("M5").Value = -1326.2
$ws.Range("N5").Value = -15524

# Row 109
$ws.Range("H109").Value = 941.6429000000001
$ws.Range("I109").Value = 254.125
$ws.Range("J109").Value = 1858.3334
$ws.Range("K109").Value = 762.375
$ws.Range("L109").Value = 5575.0002
$ws.Range("M109").Value = 277.625
$ws.Range("N109").Value = -7655.0002

# Row 122
$ws.Range("H122").Value = 1222.2142
$ws.Range("I122").Value = 559.2857
$ws.Range("J122").Value = 1885.1428
$ws.Range("K122").Value = 5033.571300000001
$ws.Range("L122").Value = 16966.2852
$ws.Range("M122").Value = -2583.571300000001
$ws.Range("N122").Value = -21866.2852

# Row 135
$ws.Range("H135").Value = 1249.5
$ws.Range("I135").Value = 479.4
$ws.Range("J135").Value = 5100
$ws.Range("K135").Value = 4314.599999999999
$ws.Range("L135").Value = 45900
$ws.Range("M135").Value = -1779.599999999999
$ws.Range("N135").Value = -50970

# Row 141
$ws.Range("H141").Value = 3075.5557
$ws.Range("J141").Value = 4772.5
$ws.Range("L141").Value = 14317.5
$ws.Range("N141").Value = -24677.5

$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 34835.332
$ws.Range("I19").Value = 4500
$ws.Range("K19").Value = 4500
$ws.Range("M19").Value = -4212

# Row 70
$ws.Range("H70").Value = 4800
$ws.Range("I70").Value = 4333.3335
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 4333.3335
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -4063.3335
$ws.Range("N70").Value = -6040

# Row 73
$ws.Range("H73").Value = 4800
$ws.Range("I73").Value = 4333.3335
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 4333.3335
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -3397.3335
$ws.Range("N73").Value = -7372

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1762.5
$ws.Range("I68").Value = 1022.7273
$ws.Range("J68").Value = 9900
$ws.Range("K68").Value = 1022.7273
$ws.Range("L68").Value = 9900
$ws.Range("M68").Value = -273.7273
$ws.Range("N68").Value = -11398

# Row 71
$ws.Range("H71").Value = 1762.5
$ws.Range("I71").Value = 1022.7273
$ws.Range("J71").Value = 9900
$ws.Range("K71").Value = 5113.636500000001
$ws.Range("L71").Value = 49500
$ws.Range("M71").Value = -1369.636500000001
$ws.Range("N71").Value = -56988

$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 65610.2
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 65610.2
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 65610.2
$ws.Range("N21").Value = -66080.2
$ws.Range("M21").ClearContents()

# Row 35
$ws.Range("H35").Value = 65610.2
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 65610.2
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 65610.2
$ws.Range("N35").Value = -66190.2
$ws.Range("M35").ClearContents()

# Row 132
$ws.Range("H132").Value = 14766.512
$ws.Range("I132").Value = 2147.9666
$ws.Range("K132").Value = 6443.899800000001
$ws.Range("M132").Value = -3913.899800000001
